# Auto-generated Excel COM-interop script
# Applies 2024-08-30 YTD update to violent-crime-ytd.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5256
$ws.Range("K3").Value = 5426
$ws.Range("K4").Value = 1125
$ws.Range("K6").Value = 6037
$ws.Range("K7").Value = 18231

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K5").Value = 42
$ws.Range("K6").Value = 133
$ws.Range("K7").Value = 538
$ws.Range("K8").Value = 1230
$ws.Range("K11").Value = 349
$ws.Range("K18").Value = 123
$ws.Range("K19").Value = 536
$ws.Range("K21").Value = 57
$ws.Range("K22").Value = 48
$ws.Range("K23").Value = 191
$ws.Range("K27").Value = 173
$ws.Range("K29").Value = 976
$ws.Range("K31").Value = 199
$ws.Range("K33").Value = 779
$ws.Range("K36").Value = 241
$ws.Range("K37").Value = 614
$ws.Range("K40").Value = 42
$ws.Range("K41").Value = 127
$ws.Range("K42").Value = 674
$ws.Range("K43").Value = 161
$ws.Range("K45").Value = 21
$ws.Range("K47").Value = 123
$ws.Range("K48").Value = 228
$ws.Range("K51").Value = 229
$ws.Range("K53").Value = 235
$ws.Range("K54").Value = 358
$ws.Range("K55").Value = 204
$ws.Range("K57").Value = 67
$ws.Range("K60").Value = 113
$ws.Range("K63").Value = 50
$ws.Range("K64").Value = 116
$ws.Range("K65").Value = 415
$ws.Range("K67").Value = 692
$ws.Range("K73").Value = 156
$ws.Range("K76").Value = 252
$ws.Range("K78").Value = 209
$ws.Range("K79").Value = 449
$ws.Range("K83").Value = 405
$ws.Range("K85").Value = 856
$ws.Range("K86").Value = 123
$ws.Range("K90").Value = 164
$ws.Range("K91").Value = 198
$ws.Range("K94").Value = 240
$ws.Range("K95").Value = 312
$ws.Range("K96").Value = 198
$ws.Range("K101").Value = 18231

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 183
$ws.Range("K7").Value = 538

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 91
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 286
$ws.Range("K7").Value = 856

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 235

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 337
$ws.Range("K6").Value = 418
$ws.Range("K7").Value = 1230

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 148
$ws.Range("K7").Value = 405

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 214
$ws.Range("K3").Value = 288
$ws.Range("K6").Value = 225
$ws.Range("K7").Value = 779

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 104
$ws.Range("K7").Value = 312

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 176
$ws.Range("K3").Value = 204
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 614

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 131
$ws.Range("K7").Value = 415

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K4").Value = 39
$ws.Range("K7").Value = 692

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 57
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 358

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 279
$ws.Range("K3").Value = 352
$ws.Range("K6").Value = 270
$ws.Range("K7").Value = 976

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 171
$ws.Range("K7").Value = 536

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 178
$ws.Range("K3").Value = 209
$ws.Range("K6").Value = 254
$ws.Range("K7").Value = 674

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 204

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 94
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 150
$ws.Range("K7").Value = 449

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 94
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 35
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 42
